$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(8)

# Move the old row 6 task ("Sending Waypoints from tablet to phone") down to
# row 7, renaming it and updating its progress/time-spent/date-complete.
$ws.Range("A7").Value = "Robot Moves Straight"
$ws.Range("B7").Value = 40115
$ws.Range("D7").Value = 0.9
$ws.Range("E7").Value = 1

# New row 6: the "Sending Waypoints from tablet to phone" task is now complete
$ws.Range("A6").Value = "Sending Waypoints from tablet to phone"
$ws.Range("B6").Value = 40111
$ws.Range("C6").Value = 40111
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1.5

# Update progress/time spent on "Create Code for moving to GPS coordinate"
$ws.Range("D5").Value = 0.9
$ws.Range("E5").Value = 5

# Update total time spent for the week
$ws.Range("F1").Value = 10.5

# Update the selected cell shown when the sheet is viewed
$ws.Activate()
$ws.Range("F2").Select()
